$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the last data row (row 6); dimension becomes A1:AH5 ---
$ws.Rows.Item(6).Delete()

# --- Widen specific data columns from 7 to 8 characters ---
$colsToWiden = @(2, 3, 7, 10, 11, 12, 13, 15, 17, 22, 24, 27, 28, 29, 30, 34)
foreach ($c in $colsToWiden) {
    $ws.Columns.Item($c).ColumnWidth = 7.166666666666667
}

# --- Replace the dataset values (rows 2-5) with the new readings ---
# Row 2 (A2:AH2)
$ws.Cells.Item(2, 1).Value = 45044.50694444445
$ws.Cells.Item(2, 2).Value = 5.237
$ws.Cells.Item(2, 3).Value = 1.607
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 2.478
$ws.Cells.Item(2, 6).Value = 3.053
$ws.Cells.Item(2, 7).Value = 2.519
$ws.Cells.Item(2, 8).Value = 5.488
$ws.Cells.Item(2, 9).Value = 1.554
$ws.Cells.Item(2, 10).Value = 0.9409999999999999
$ws.Cells.Item(2, 11).Value = 4.022
$ws.Cells.Item(2, 12).Value = 1.069
$ws.Cells.Item(2, 13).Value = 0.9379999999999999
$ws.Cells.Item(2, 14).Value = 0.6929999999999999
$ws.Cells.Item(2, 15).Value = 0.87
$ws.Cells.Item(2, 16).Value = 2.749
$ws.Cells.Item(2, 17).Value = 1.106
$ws.Cells.Item(2, 18).Value = 0.51
$ws.Cells.Item(2, 19).Value = 0.063
$ws.Cells.Item(2, 20).Value = 20.094
$ws.Cells.Item(2, 21).Value = 4.803
$ws.Cells.Item(2, 22).Value = 2.372
$ws.Cells.Item(2, 23).Value = 3.891
$ws.Cells.Item(2, 24).Value = 1.034
$ws.Cells.Item(2, 25).Value = 0.249
$ws.Cells.Item(2, 26).Value = 1.98
$ws.Cells.Item(2, 27).Value = 1.154
$ws.Cells.Item(2, 28).Value = 0.674
$ws.Cells.Item(2, 29).Value = 0.9399999999999999
$ws.Cells.Item(2, 30).Value = 3.129
$ws.Cells.Item(2, 31).Value = 2.866
$ws.Cells.Item(2, 32).Value = 3.128
$ws.Cells.Item(2, 33).Value = 0.446
$ws.Cells.Item(2, 34).Value = 1.43

# Row 3 (A3:AH3)
$ws.Cells.Item(3, 1).Value = 45044.51388888889
$ws.Cells.Item(3, 2).Value = 20.388
$ws.Cells.Item(3, 3).Value = 14.533
$ws.Cells.Item(3, 4).Value = 0.481
$ws.Cells.Item(3, 5).Value = 40.466
$ws.Cells.Item(3, 6).Value = 34.136
$ws.Cells.Item(3, 7).Value = 15.603
$ws.Cells.Item(3, 8).Value = 51.591
$ws.Cells.Item(3, 9).Value = 22.572
$ws.Cells.Item(3, 10).Value = 10.373
$ws.Cells.Item(3, 11).Value = 16.498
$ws.Cells.Item(3, 12).Value = 16.371
$ws.Cells.Item(3, 13).Value = 17.046
$ws.Cells.Item(3, 14).Value = 4.844
$ws.Cells.Item(3, 15).Value = 14.532
$ws.Cells.Item(3, 16).Value = 21.861
$ws.Cells.Item(3, 17).Value = 12.194
$ws.Cells.Item(3, 18).Value = 0.533
$ws.Cells.Item(3, 19).Value = 0.474
$ws.Cells.Item(3, 20).Value = 221.152
$ws.Cells.Item(3, 21).Value = 41.875
$ws.Cells.Item(3, 22).Value = 14.215
$ws.Cells.Item(3, 23).Value = 28.576
$ws.Cells.Item(3, 24).Value = 14.582
$ws.Cells.Item(3, 25).Value = 1.972
$ws.Cells.Item(3, 26).Value = 26.371
$ws.Cells.Item(3, 27).Value = 12.17
$ws.Cells.Item(3, 28).Value = 10.413
$ws.Cells.Item(3, 29).Value = 12.358
$ws.Cells.Item(3, 30).Value = 18.238
$ws.Cells.Item(3, 31).Value = 1.133
$ws.Cells.Item(3, 32).Value = 46.087
$ws.Cells.Item(3, 33).Value = 7.62
$ws.Cells.Item(3, 34).Value = 17.006

# Row 4 (A4:AH4)
$ws.Cells.Item(4, 1).Value = 45044.52083333334
$ws.Cells.Item(4, 2).Value = 25.239
$ws.Cells.Item(4, 3).Value = 18.517
$ws.Cells.Item(4, 4).Value = 0.679
$ws.Cells.Item(4, 5).Value = 52.414
$ws.Cells.Item(4, 6).Value = 43.813
$ws.Cells.Item(4, 7).Value = 19.638
$ws.Cells.Item(4, 8).Value = 74.30800000000001
$ws.Cells.Item(4, 9).Value = 29.178
$ws.Cells.Item(4, 10).Value = 13.306
$ws.Cells.Item(4, 11).Value = 20.459
$ws.Cells.Item(4, 12).Value = 21.163
$ws.Cells.Item(4, 13).Value = 22.132
$ws.Cells.Item(4, 14).Value = 6.165
$ws.Cells.Item(4, 15).Value = 18.845
$ws.Cells.Item(4, 16).Value = 27.747
$ws.Cells.Item(4, 17).Value = 15.694
$ws.Cells.Item(4, 18).Value = 0.495
$ws.Cells.Item(4, 19).Value = 0.63
$ws.Cells.Item(4, 20).Value = 284.29
$ws.Cells.Item(4, 21).Value = 53.615
$ws.Cells.Item(4, 22).Value = 17.921
$ws.Cells.Item(4, 23).Value = 36.418
$ws.Cells.Item(4, 24).Value = 18.829
$ws.Cells.Item(4, 25).Value = 2.514
$ws.Cells.Item(4, 26).Value = 36.494
$ws.Cells.Item(4, 27).Value = 15.593
$ws.Cells.Item(4, 28).Value = 13.485
$ws.Cells.Item(4, 29).Value = 15.939
$ws.Cells.Item(4, 30).Value = 22.906
$ws.Cells.Item(4, 31).Value = 0.716
$ws.Cells.Item(4, 32).Value = 67.23099999999999
$ws.Cells.Item(4, 33).Value = 9.890000000000001
$ws.Cells.Item(4, 34).Value = 21.898

# Row 5 (A5:AH5)
$ws.Cells.Item(5, 1).Value = 45044.52777777778
$ws.Cells.Item(5, 2).Value = 14.47
$ws.Cells.Item(5, 3).Value = 10.57
$ws.Cells.Item(5, 4).Value = 0.35
$ws.Cells.Item(5, 5).Value = 29.64
$ws.Cells.Item(5, 6).Value = 24.88
$ws.Cells.Item(5, 7).Value = 11.25
$ws.Cells.Item(5, 8).Value = 46.03
$ws.Cells.Item(5, 9).Value = 16.49
$ws.Cells.Item(5, 10).Value = 7.58
$ws.Cells.Item(5, 11).Value = 11.72
$ws.Cells.Item(5, 12).Value = 12.02
$ws.Cells.Item(5, 13).Value = 12.48
$ws.Cells.Item(5, 14).Value = 3.51
$ws.Cells.Item(5, 15).Value = 10.66
$ws.Cells.Item(5, 16).Value = 15.88
$ws.Cells.Item(5, 17).Value = 8.890000000000001
$ws.Cells.Item(5, 18).Value = 0.35
$ws.Cells.Item(5, 19).Value = 0.33
$ws.Cells.Item(5, 20).Value = 158.46
$ws.Cells.Item(5, 21).Value = 30.58
$ws.Cells.Item(5, 22).Value = 10.23
$ws.Cells.Item(5, 23).Value = 20.77
$ws.Cells.Item(5, 24).Value = 10.67
$ws.Cells.Item(5, 25).Value = 1.43
$ws.Cells.Item(5, 26).Value = 21.96
$ws.Cells.Item(5, 27).Value = 8.869999999999999
$ws.Cells.Item(5, 28).Value = 7.63
$ws.Cells.Item(5, 29).Value = 9.029999999999999
$ws.Cells.Item(5, 30).Value = 13.11
$ws.Cells.Item(5, 31).Value = 0.52
$ws.Cells.Item(5, 32).Value = 41.63
$ws.Cells.Item(5, 33).Value = 5.58
$ws.Cells.Item(5, 34).Value = 12.41
